$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Convert E38 and E39 from text to real numbers (matching the diff) ---
$ws.Cells.Item(38, 5).Value = 20
$ws.Cells.Item(39, 5).Value = 531344

# --- Row 40 ---
$ws.Cells.Item(40, 1).Value = "24/06/2024 08:44:40"
$ws.Cells.Item(40, 2).Value = 1
$ws.Cells.Item(40, 3).Value = "BSE"
$ws.Cells.Item(40, 4).Value = "BSE (Bombay stock exchange)"

$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "20"
$ws.Cells.Item(40, 5).Style = "Normal"

$ws.Cells.Item(40, 6).Value = -2.37
$ws.Cells.Item(40, 7).Value = 2499.6
$ws.Cells.Item(40, 8).Value = 754277

# --- Row 41 ---
$ws.Cells.Item(41, 1).Value = "24/06/2024 08:44:40"
$ws.Cells.Item(41, 2).Value = 2
$ws.Cells.Item(41, 3).Value = "CONCOR"
$ws.Cells.Item(41, 4).Value = "Container Corporation Of India Limited"

$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "531344"
$ws.Cells.Item(41, 5).Style = "Normal"

$ws.Cells.Item(41, 6).Value = -3.38
$ws.Cells.Item(41, 7).Value = 1054.05
$ws.Cells.Item(41, 8).Value = 4044432
